$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F3 previously shared the "Comment" header string; give it its own text.
$ws.Range("F3").Value = "Hello world"

# F4 previously held "Available comment"; update its text.
$ws.Range("F4").Value = "Hello, world"
